$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Silent (沉默)
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "沉默"
$ws.Cells.Item(14, 3).Value = "Silent"
$ws.Cells.Item(14, 4).Value = -1
$ws.Cells.Item(14, 5).Value = -1
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = "buff_silent_icon"
$ws.Cells.Item(14, 11).Value = "沉默的单位无法进行攻击"

# Row 15: Confine (禁锢)
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "禁锢"
$ws.Cells.Item(15, 3).Value = "Confine"
$ws.Cells.Item(15, 4).Value = -1
$ws.Cells.Item(15, 5).Value = -1
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = "buff_confine_icon"
$ws.Cells.Item(15, 11).Value = "禁锢的单位无法移动"

$ws.Range("K12").Select()
